$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item("TextBox 12")
$tr = $shp.TextFrame.TextRange
$tr.Delete()
[void]$tr.InsertAfter("Comfortable Death")
